# Adding pickling to replication
# Update the bootstrap standard-error values in columns C and D for the
# "year range" rows on Sheet1 to the newly re-computed (pickled) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "(0.01)"
$ws.Range("D3").Value = "(0.01)"

$ws.Range("C5").Value = "(0.71)"
$ws.Range("D5").Value = "(0.49)"

$ws.Range("C7").Value = "(0.2)"
$ws.Range("D7").Value = "(0.02)"

$ws.Range("C9").Value = "(0.48)"
$ws.Range("D9").Value = "(0.65)"

$ws.Range("C11").Value = "(0.46)"
$ws.Range("D11").Value = "(0.48)"

$ws.Range("C13").Value = "(0.02)"
$ws.Range("D13").Value = "(0.03)"

$ws.Range("C15").Value = "(0.42)"
$ws.Range("D15").Value = "(0.4)"
